$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 5).Value = 2
$ws.Cells.Item(2, 6).Value = 0.6666666666666666
$ws.Cells.Item(2, 7).Value = 0.1751946666666667
$ws.Cells.Item(2, 8).Value = 0.5255840000000001
$ws.Cells.Item(2, 9).Value = 0.07933130249481599
$ws.Cells.Item(2, 10).Value = 0.079331302494816
$ws.Cells.Item(2, 11).Value = 1
$ws.Cells.Item(2, 12).Value = 0.3333333333333333
$ws.Cells.Item(2, 13).Value = 0.04996866666666667
$ws.Cells.Item(2, 14).Value = 0.149906
$ws.Cells.Item(2, 15).Value = 0.06831051926220302
$ws.Cells.Item(2, 16).Value = 0.06831051926220301
$ws.Cells.Item(2, 17).Value = 0.008754243900444447
$ws.Cells.Item(2, 18).Value = 0.07878819510400001
$ws.Cells.Item(2, 19).Value = 0.005419162467167782
$ws.Cells.Item(2, 20).Value = 0.005419162467167782

# Row 3
$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(3, 6).Value = 0.6666666666666666
$ws.Cells.Item(3, 7).Value = 0.1751946666666667
$ws.Cells.Item(3, 8).Value = 0.5255840000000001
$ws.Cells.Item(3, 9).Value = 0.07933130249481599
$ws.Cells.Item(3, 10).Value = 0.079331302494816
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 0.5898666666666667
$ws.Cells.Item(3, 14).Value = 1.7696
$ws.Cells.Item(3, 15).Value = 0.8063873019518528
$ws.Cells.Item(3, 16).Value = 0.8063873019518527
$ws.Cells.Item(3, 17).Value = 0.1033414940444445
$ws.Cells.Item(3, 18).Value = 0.9300734464000001
$ws.Cells.Item(3, 19).Value = 0.06397175497912096
$ws.Cells.Item(3, 20).Value = 0.06397175497912096

# Row 4
$ws.Cells.Item(4, 5).Value = 2
$ws.Cells.Item(4, 6).Value = 0.6666666666666666
$ws.Cells.Item(4, 7).Value = 0.1751946666666667
$ws.Cells.Item(4, 8).Value = 0.5255840000000001
$ws.Cells.Item(4, 9).Value = 0.07933130249481599
$ws.Cells.Item(4, 10).Value = 0.079331302494816
$ws.Cells.Item(4, 11).Value = 2
$ws.Cells.Item(4, 12).Value = 0.6666666666666666
$ws.Cells.Item(4, 13).Value = 0.09165766666666668
$ws.Cells.Item(4, 14).Value = 0.274973
$ws.Cells.Item(4, 15).Value = 0.1253021787859442
$ws.Cells.Item(4, 16).Value = 0.1253021787859442
$ws.Cells.Item(4, 17).Value = 0.01605793435911112
$ws.Cells.Item(4, 18).Value = 0.144521409232
$ws.Cells.Item(4, 19).Value = 0.009940385048527255
$ws.Cells.Item(4, 20).Value = 0.009940385048527255

# Row 5
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 0.8048609999999999
$ws.Cells.Item(5, 8).Value = 2.414583
$ws.Cells.Item(5, 9).Value = 0.364455566325916
$ws.Cells.Item(5, 10).Value = 0.3644555663259161
$ws.Cells.Item(5, 11).Value = 1
$ws.Cells.Item(5, 12).Value = 0.3333333333333333
$ws.Cells.Item(5, 13).Value = 0.04996866666666667
$ws.Cells.Item(5, 14).Value = 0.149906
$ws.Cells.Item(5, 15).Value = 0.06831051926220302
$ws.Cells.Item(5, 16).Value = 0.06831051926220301
$ws.Cells.Item(5, 17).Value = 0.040217831022
$ws.Cells.Item(5, 18).Value = 0.361960479198
$ws.Cells.Item(5, 19).Value = 0.0248961489837236
$ws.Cells.Item(5, 20).Value = 0.0248961489837236

# Row 6
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 0.8048609999999999
$ws.Cells.Item(6, 8).Value = 2.414583
$ws.Cells.Item(6, 9).Value = 0.364455566325916
$ws.Cells.Item(6, 10).Value = 0.3644555663259161
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 0.5898666666666667
$ws.Cells.Item(6, 14).Value = 1.7696
$ws.Cells.Item(6, 15).Value = 0.8063873019518528
$ws.Cells.Item(6, 16).Value = 0.8063873019518527
$ws.Cells.Item(6, 17).Value = 0.4747606751999999
$ws.Cells.Item(6, 18).Value = 4.2728460768
$ws.Cells.Item(6, 19).Value = 0.2938923408108899
$ws.Cells.Item(6, 20).Value = 0.29389234081089

# Row 7
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 0.8048609999999999
$ws.Cells.Item(7, 8).Value = 2.414583
$ws.Cells.Item(7, 9).Value = 0.364455566325916
$ws.Cells.Item(7, 10).Value = 0.3644555663259161
$ws.Cells.Item(7, 11).Value = 2
$ws.Cells.Item(7, 12).Value = 0.6666666666666666
$ws.Cells.Item(7, 13).Value = 0.09165766666666668
$ws.Cells.Item(7, 14).Value = 0.274973
$ws.Cells.Item(7, 15).Value = 0.1253021787859442
$ws.Cells.Item(7, 16).Value = 0.1253021787859442
$ws.Cells.Item(7, 17).Value = 0.073771681251
$ws.Cells.Item(7, 18).Value = 0.663945131259
$ws.Cells.Item(7, 19).Value = 0.04566707653130248
$ws.Cells.Item(7, 20).Value = 0.04566707653130248

# Row 8
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 1.142989
$ws.Cells.Item(8, 8).Value = 3.428967
$ws.Cells.Item(8, 9).Value = 0.5175660186035755
$ws.Cells.Item(8, 10).Value = 0.5175660186035757
$ws.Cells.Item(8, 11).Value = 1
$ws.Cells.Item(8, 12).Value = 0.3333333333333333
$ws.Cells.Item(8, 13).Value = 0.04996866666666667
$ws.Cells.Item(8, 14).Value = 0.149906
$ws.Cells.Item(8, 15).Value = 0.06831051926220302
$ws.Cells.Item(8, 16).Value = 0.06831051926220301
$ws.Cells.Item(8, 17).Value = 0.05711363634466667
$ws.Cells.Item(8, 18).Value = 0.514022727102
$ws.Cells.Item(8, 19).Value = 0.03535520348328128
$ws.Cells.Item(8, 20).Value = 0.03535520348328128

# Row 9
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 1.142989
$ws.Cells.Item(9, 8).Value = 3.428967
$ws.Cells.Item(9, 9).Value = 0.5175660186035755
$ws.Cells.Item(9, 10).Value = 0.5175660186035757
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 0.5898666666666667
$ws.Cells.Item(9, 14).Value = 1.7696
$ws.Cells.Item(9, 15).Value = 0.8063873019518528
$ws.Cells.Item(9, 16).Value = 0.8063873019518527
$ws.Cells.Item(9, 17).Value = 0.6742111114666667
$ws.Cells.Item(9, 18).Value = 6.0679000032
$ws.Cells.Item(9, 19).Value = 0.4173586653236998
$ws.Cells.Item(9, 20).Value = 0.4173586653236998

# Row 10
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 1.142989
$ws.Cells.Item(10, 8).Value = 3.428967
$ws.Cells.Item(10, 9).Value = 0.5175660186035755
$ws.Cells.Item(10, 10).Value = 0.5175660186035757
$ws.Cells.Item(10, 11).Value = 2
$ws.Cells.Item(10, 12).Value = 0.6666666666666666
$ws.Cells.Item(10, 13).Value = 0.09165766666666668
$ws.Cells.Item(10, 14).Value = 0.274973
$ws.Cells.Item(10, 15).Value = 0.1253021787859442
$ws.Cells.Item(10, 16).Value = 0.1253021787859442
$ws.Cells.Item(10, 17).Value = 0.1047637047656667
$ws.Cells.Item(10, 18).Value = 0.9428733428910001
$ws.Cells.Item(10, 19).Value = 0.06485214979659455
$ws.Cells.Item(10, 20).Value = 0.06485214979659455

# Row 11
$ws.Cells.Item(11, 5).Value = 2
$ws.Cells.Item(11, 6).Value = 0.6666666666666666
$ws.Cells.Item(11, 7).Value = 0.08534799999999999
$ws.Cells.Item(11, 8).Value = 0.256044
$ws.Cells.Item(11, 9).Value = 0.0386471125756923
$ws.Cells.Item(11, 10).Value = 0.0386471125756923
$ws.Cells.Item(11, 11).Value = 1
$ws.Cells.Item(11, 12).Value = 0.3333333333333333
$ws.Cells.Item(11, 13).Value = 0.04996866666666667
$ws.Cells.Item(11, 14).Value = 0.149906
$ws.Cells.Item(11, 15).Value = 0.06831051926220302
$ws.Cells.Item(11, 16).Value = 0.06831051926220301
$ws.Cells.Item(11, 17).Value = 0.004264725762666666
$ws.Cells.Item(11, 18).Value = 0.038382531864
$ws.Cells.Item(11, 19).Value = 0.002640004328030357
$ws.Cells.Item(11, 20).Value = 0.002640004328030357

# Row 12
$ws.Cells.Item(12, 5).Value = 2
$ws.Cells.Item(12, 6).Value = 0.6666666666666666
$ws.Cells.Item(12, 7).Value = 0.08534799999999999
$ws.Cells.Item(12, 8).Value = 0.256044
$ws.Cells.Item(12, 9).Value = 0.0386471125756923
$ws.Cells.Item(12, 10).Value = 0.0386471125756923
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 12).Value = 1
$ws.Cells.Item(12, 13).Value = 0.5898666666666667
$ws.Cells.Item(12, 14).Value = 1.7696
$ws.Cells.Item(12, 15).Value = 0.8063873019518528
$ws.Cells.Item(12, 16).Value = 0.8063873019518527
$ws.Cells.Item(12, 17).Value = 0.05034394026666666
$ws.Cells.Item(12, 18).Value = 0.4530954624
$ws.Cells.Item(12, 19).Value = 0.03116454083814203
$ws.Cells.Item(12, 20).Value = 0.03116454083814203

# Row 13
$ws.Cells.Item(13, 5).Value = 2
$ws.Cells.Item(13, 6).Value = 0.6666666666666666
$ws.Cells.Item(13, 7).Value = 0.08534799999999999
$ws.Cells.Item(13, 8).Value = 0.256044
$ws.Cells.Item(13, 9).Value = 0.0386471125756923
$ws.Cells.Item(13, 10).Value = 0.0386471125756923
$ws.Cells.Item(13, 11).Value = 2
$ws.Cells.Item(13, 12).Value = 0.6666666666666666
$ws.Cells.Item(13, 13).Value = 0.09165766666666668
$ws.Cells.Item(13, 14).Value = 0.274973
$ws.Cells.Item(13, 15).Value = 0.1253021787859442
$ws.Cells.Item(13, 16).Value = 0.1253021787859442
$ws.Cells.Item(13, 17).Value = 0.007822798534666667
$ws.Cells.Item(13, 18).Value = 0.070405186812
$ws.Cells.Item(13, 19).Value = 0.004842567409519909
$ws.Cells.Item(13, 20).Value = 0.004842567409519909

